# Applies scheduled market-data refresh updates to the Sheets workbook.
# Each worksheet corresponds to a crafting class (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For every affected leve row, the pricing/profit columns (H-N) are updated to the
# latest computed values. Where a cell no longer has data it is cleared; where a
# cell gains data for the first time it is set.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 12
$ws.Range("H12").Value = 25272
$ws.Range("I12").Value = 25272
$ws.Range("K12").Value = 25272
$ws.Range("M12").Value = -25102

# Row 19
$ws.Range("H19").Value = 1440
$ws.Range("I19").Value = 1190
$ws.Range("K19").Value = 1190
$ws.Range("M19").Value = -1015

# Row 40
$ws.Range("H40").Value = 6813.143
$ws.Range("J40").Value = 10332.333
$ws.Range("L40").Value = 10332.333
$ws.Range("N40").Value = -10682.333

# Row 74
$ws.Range("H74").Value = 1806758
$ws.Range("I74").Value = 1806758
$ws.Range("K74").Value = 1806758
$ws.Range("M74").Value = -1805822

# Row 76
$ws.Range("H76").Value = 3928.4285
$ws.Range("I76").Value = 3683
$ws.Range("K76").Value = 3683
$ws.Range("M76").Value = -3368

# Row 77
$ws.Range("H77").Value = 1806758
$ws.Range("I77").Value = 1806758
$ws.Range("K77").Value = 9033790
$ws.Range("M77").Value = -9029110

# Row 79
$ws.Range("H79").Value = 3928.4285
$ws.Range("I79").Value = 3683
$ws.Range("K79").Value = 3683
$ws.Range("M79").Value = -2591

# Row 80
$ws.Range("H80").Value = 2507
$ws.Range("J80").Value = 587.5
$ws.Range("L80").Value = 1762.5
$ws.Range("N80").Value = -3758.5

# Row 83
$ws.Range("H83").Value = 2507
$ws.Range("J83").Value = 587.5
$ws.Range("L83").Value = 5287.5
$ws.Range("N83").Value = -15271.5

# Row 88
$ws.Range("H88").Value = 21720860
$ws.Range("I88").Value = 66672668
$ws.Range("J88").Value = 2990939
$ws.Range("K88").Value = 66672668
$ws.Range("L88").Value = 2990939
$ws.Range("M88").Value = -66672262
$ws.Range("N88").Value = -2991751

# Row 91
$ws.Range("H91").Value = 21720860
$ws.Range("I91").Value = 66672668
$ws.Range("J91").Value = 2990939
$ws.Range("K91").Value = 66672668
$ws.Range("L91").Value = 2990939
$ws.Range("M91").Value = -66671264
$ws.Range("N91").Value = -2993747

# Row 116
$ws.Range("H116").Value = 4700.5
$ws.Range("I116").Value = 4875.75
$ws.Range("J116").Value = 3999.5
$ws.Range("K116").Value = 4875.75
$ws.Range("L116").Value = 3999.5
$ws.Range("M116").Value = -1433.75
$ws.Range("N116").Value = -10883.5

$ws = $wb.Worksheets.Item("ARM")

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""

# Row 45
$ws.Range("H45").Value = 8703.817999999999
$ws.Range("I45").Value = 10114.375
$ws.Range("J45").Value = 4942.3335
$ws.Range("K45").Value = 10114.375
$ws.Range("L45").Value = 4942.3335
$ws.Range("M45").Value = -9737.375
$ws.Range("N45").Value = -5696.3335

# Row 63
$ws.Range("H63").Value = 1759.909
$ws.Range("I63").Value = 1795
$ws.Range("K63").Value = 1795
$ws.Range("M63").Value = -1109

# Row 66
$ws.Range("H66").Value = 1759.909
$ws.Range("I66").Value = 1795
$ws.Range("K66").Value = 8975
$ws.Range("M66").Value = -5543

$ws = $wb.Worksheets.Item("BSM")

# Row 36
$ws.Range("H36").Value = 7299.6665
$ws.Range("I36").Value = 7299.6665
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 7299.6665
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -6765.6665
$ws.Range("N36").Value = ""

# Row 39
$ws.Range("H39").Value = 33499.5
$ws.Range("I39").Value = 32000
$ws.Range("K39").Value = 32000
$ws.Range("M39").Value = -31611

# Row 64
$ws.Range("H64").Value = 759.6667
$ws.Range("I64").Value = 623
$ws.Range("J64").Value = 828
$ws.Range("K64").Value = 623
$ws.Range("L64").Value = 828
$ws.Range("M64").Value = -398
$ws.Range("N64").Value = -1278

# Row 67
$ws.Range("H67").Value = 759.6667
$ws.Range("I67").Value = 623
$ws.Range("J67").Value = 828
$ws.Range("K67").Value = 623
$ws.Range("L67").Value = 828
$ws.Range("M67").Value = 157
$ws.Range("N67").Value = -2388

# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""

# Row 134
$ws.Range("H134").Value = 17247072
$ws.Range("I134").Value = 20006130
$ws.Range("J134").Value = 2965.25
$ws.Range("K134").Value = 60018390
$ws.Range("L134").Value = 8895.75
$ws.Range("M134").Value = -60015855
$ws.Range("N134").Value = -13965.75

$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 1587.4546
$ws.Range("I16").Value = 1657.75
$ws.Range("J16").Value = 1400
$ws.Range("K16").Value = 1657.75
$ws.Range("L16").Value = 1400
$ws.Range("M16").Value = -1370.75
$ws.Range("N16").Value = -1974

# Row 31
$ws.Range("H31").Value = 4821.84
$ws.Range("I31").Value = 3302.4211
$ws.Range("K31").Value = 3302.4211
$ws.Range("M31").Value = -3007.4211

# Row 34
$ws.Range("H34").Value = 4821.84
$ws.Range("I34").Value = 3302.4211
$ws.Range("K34").Value = 3302.4211
$ws.Range("M34").Value = -3100.4211

# Row 113
$ws.Range("H113").Value = 1587.4546
$ws.Range("I113").Value = 1657.75
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 1657.75
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = 512.25
$ws.Range("N113").Value = -5740

# Row 122
$ws.Range("H122").Value = 1842.125
$ws.Range("I122").Value = 2112.4546
$ws.Range("J122").Value = 1247.4
$ws.Range("K122").Value = 6337.3638
$ws.Range("L122").Value = 3742.2
$ws.Range("M122").Value = -3887.3638
$ws.Range("N122").Value = -8642.200000000001

$ws = $wb.Worksheets.Item("CUL")

# Row 11
$ws.Range("H11").Value = 189687.5
$ws.Range("J11").Value = 100000
$ws.Range("L11").Value = 300000
$ws.Range("N11").Value = -300280

# Row 64
$ws.Range("H64").Value = 6124.5
$ws.Range("J64").Value = 10999
$ws.Range("L64").Value = 32997
$ws.Range("N64").Value = -33537

# Row 67
$ws.Range("H67").Value = 6124.5
$ws.Range("J67").Value = 10999
$ws.Range("L67").Value = 32997
$ws.Range("N67").Value = -34869

# Row 68
$ws.Range("H68").Value = 2461.9443
$ws.Range("I68").Value = 2646.25
$ws.Range("J68").Value = 2314.5
$ws.Range("K68").Value = 7938.75
$ws.Range("L68").Value = 6943.5
$ws.Range("M68").Value = -7127.75
$ws.Range("N68").Value = -8565.5

# Row 71
$ws.Range("H71").Value = 2461.9443
$ws.Range("I71").Value = 2646.25
$ws.Range("J71").Value = 2314.5
$ws.Range("K71").Value = 23816.25
$ws.Range("L71").Value = 20830.5
$ws.Range("M71").Value = -19760.25
$ws.Range("N71").Value = -28942.5

# Row 129
$ws.Range("H129").Value = 2824.65
$ws.Range("I129").Value = 720.4286
$ws.Range("J129").Value = 3957.6924
$ws.Range("K129").Value = 2161.2858
$ws.Range("L129").Value = 11873.0772
$ws.Range("M129").Value = 2838.7142
$ws.Range("N129").Value = -21873.0772

# Row 131
$ws.Range("H131").Value = 2044.7142
$ws.Range("I131").Value = 2110
$ws.Range("K131").Value = 6330
$ws.Range("M131").Value = -1290

# Row 139
$ws.Range("H139").Value = 7000
$ws.Range("I139").Value = 7000
$ws.Range("K139").Value = 21000
$ws.Range("M139").Value = -15860

$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 6100.5625
$ws.Range("I70").Value = 5900.846
$ws.Range("K70").Value = 5900.846
$ws.Range("M70").Value = -5630.846

# Row 73
$ws.Range("H73").Value = 6100.5625
$ws.Range("I73").Value = 5900.846
$ws.Range("K73").Value = 5900.846
$ws.Range("M73").Value = -4964.846

# Row 80
$ws.Range("H80").Value = 4020
$ws.Range("I80").Value = 4223
$ws.Range("K80").Value = 4223
$ws.Range("M80").Value = -3225

# Row 83
$ws.Range("H83").Value = 4020
$ws.Range("I83").Value = 4223
$ws.Range("K83").Value = 21115
$ws.Range("M83").Value = -16123

# Row 92
$ws.Range("H92").Value = 13263.667
$ws.Range("J92").Value = 13263.667
$ws.Range("L92").Value = 13263.667
$ws.Range("N92").Value = -17007.667

# Row 97
$ws.Range("H97").Value = 1647
$ws.Range("J97").Value = 1825
$ws.Range("L97").Value = 1825
$ws.Range("N97").Value = -2817

# Row 122
$ws.Range("H122").Value = 7511.25
$ws.Range("I122").Value = 5348.5
$ws.Range("K122").Value = 16045.5
$ws.Range("M122").Value = -13595.5

$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 3251.5
$ws.Range("I22").Value = 3168.8333
$ws.Range("K22").Value = 3168.8333
$ws.Range("M22").Value = -2873.8333

# Row 27
$ws.Range("H27").Value = 3251.5
$ws.Range("I27").Value = 3168.8333
$ws.Range("K27").Value = 3168.8333
$ws.Range("M27").Value = -3061.8333

# Row 40
$ws.Range("H40").Value = 3380.2144
$ws.Range("I40").Value = 3380.2144
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3380.2144
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3244.2144
$ws.Range("N40").Value = ""

# Row 82
$ws.Range("H82").Value = 1614.3846
$ws.Range("I82").Value = 1453.5454
$ws.Range("J82").Value = 2499
$ws.Range("K82").Value = 1453.5454
$ws.Range("L82").Value = 2499
$ws.Range("M82").Value = -1092.5454
$ws.Range("N82").Value = -3221

# Row 85
$ws.Range("H85").Value = 1614.3846
$ws.Range("I85").Value = 1453.5454
$ws.Range("J85").Value = 2499
$ws.Range("K85").Value = 1453.5454
$ws.Range("L85").Value = 2499
$ws.Range("M85").Value = -205.5454
$ws.Range("N85").Value = -4995

$ws = $wb.Worksheets.Item("WVR")

# Row 123
$ws.Range("H123").Value = 64999.5
$ws.Range("J123").Value = 64999.5
$ws.Range("L123").Value = 64999.5
$ws.Range("N123").Value = -74799.5

# Row 125
$ws.Range("H125").Value = 63750
$ws.Range("J125").Value = 63750
$ws.Range("L125").Value = 63750
$ws.Range("N125").Value = -73590

# Row 126
$ws.Range("H126").Value = 2793.4211
$ws.Range("I126").Value = 2716.1765
$ws.Range("K126").Value = 8148.529500000001
$ws.Range("M126").Value = -5678.529500000001
